$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Total" row values to reflect the newly purchased stocks
# (BHEL, TATA CONSUMER, RELIANCE)
$ws.Range("B2").Value = 246169.15
$ws.Range("C2").Value = 265882.5480957031
$ws.Range("D2").Value = 19713.39809570313
$ws.Range("E2").Value = 8.00807009964617
